# Applies the zh (Traditional/Simplified-source) -> English translation
# edits described by the commit diff, using Find/Replace (ReplaceAll)
# over the whole document content.

$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Language-switcher line at top of doc (occurs twice, same replacement both
# times: the hyperlinked "英語" label and the plain-text "英語" heading below).
Replace-All "英語" "English"

# Language-switcher separator list. (Search text intentionally omits the
# leading space/slash, which sit right at the hyperlink run boundary —
# including them causes Find/Replace to inherit the hyperlink's character
# formatting instead of keeping the host run's own rPr.)
Replace-All "葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語" "Portuguese / French / Thai / Vietnamese / Spanish"

# Table: "Brief" label + description.
Replace-All "簡介" "Brief"
Replace-All "發送給目標國家中那些文件未通過我們驗證流程的合作夥伴的電子郵件。 將通過 customer.io 發送" "An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io"

# Table: "Target audience" label + description.
Replace-All "目標受眾" "Target audience"
Replace-All "提交錯誤/不完整文檔的被邀請合作夥伴" "Invited partners who submitted wrong/incomplete documents"

# Subject line.
Replace-All "主題行" "Subject line"
Replace-All "[事件名稱]" "[EVENT NAME]"
Replace-All " — 文件驗證失敗 " " — document verification failed "

# Headline.
Replace-All "啊哦！ 文檔無法驗證" "Uh oh! Your documents couldn’t be verified"

# Greeting.
Replace-All "您好 " "Hi "
Replace-All "[合作夥伴姓名]" "[PARTNER NAME]"

# Intro sentence.
Replace-All "很遺憾地通知您，您的文檔未通過驗證流程，因為我們發現以下問題： " "We regret to inform you that your documents have failed our verification process as we found the following issues with them: "

# Bulleted issue list.
Replace-All "您的疫苗接種證明副本" "A copy of your vaccination certificate"
Replace-All ": 文檔不清楚" ": Document is unclear"
Replace-All "[文檔 2]" "[Document 2]"
Replace-All ": [問題]" ": [problem]"

# Resubmission deadline sentence.
Replace-All "請在 " "Please resubmit the documents above by "
Replace-All "日月年" "DD Mmm YYYY"
Replace-All " 之前重新提交上述文檔，以便我們進行必要的安排。" " so we can proceed with the necessary arrangements."

# Closing contact paragraph (must run "發送..." replacement above before this,
# since "。 " is also a substring of that longer run).
Replace-All "如有任何疑問，請通過 " "If you have any questions, please contact your country manager, "
Replace-All "[電子郵件地址]" "[NAME]"
Replace-All " 或 " ", at "
Replace-All "[WHATSAPP 號碼]" "[EMAIL ADDRESS]"
Replace-All " (WhatsApp) 聯繫您的區域經理, " " or "
Replace-All "[姓名]" "[WHATSAPP NO]"
Replace-All "。 " " (WhatsApp). "
